# Reshapes the "limites_lrf" sheet from a metric-per-row / period-per-column
# layout into a period-per-row / metric-per-column layout.
#
# Before:
#   row1:      B1:F1  = period names
#   rows3-12:  A col  = metric name, B:F = value for each period
#
# After:
#   row1:      B1:K1  = metric names (header row)
#   rows2-6:   A col  = period name, B:K = value for each metric

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - the whole used range is being replaced.
$ws.Cells.Clear()

$numFmt = "#,##0.00"
$dash   = " -   "

# ---- Header row (metric names across columns B..K) ----------------------
$headers = @(
  "RECEITA CORRENTE LÍQUIDA – RCL (IV) ",
  "(-) Transferências Obrigatórias da União relativas às emendas individuais (V) (§ 1°, art. 166 da CF) ",
  "(-) Transferências obrigatórias da União relativas às emendas de bancada (art. 166, § 16 da CF) (VI)",
  "(-) Transferências da União relativas à remuneração dos agentes comunitários de saúde e de combate às endemias (CF. art.498, § 11)",
  "(-) Outras Deduções Constitucionais ou Legais",
  "RECEITA CORRENTE LÍQUIDA AJUSTADA PARA CÁLCULO DOS LIMITES DA DESPESA COM PESSOAL (VII) = (IV - V - VI)",
  "DESPESA TOTAL COM PESSOAL – DTP sobre a RCL (VII) = (III a + III b)",
  "LIMITE MÁXIMO (IX) (incisos I, II e III art. 20 da LRF) – 49,00%",
  "LIMITE PRUDENCIAL (X) (parágrafo único, art. 22 da LRF) – 46,55%",
  "LIMITE DE ALERTA (XI) (inciso II do § 1º do art. 59 da LRF) – 44,10%"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $col = $i + 2   # B=2 .. K=11
  $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# ---- Data rows (one per period) ------------------------------------------
# Each row: period label, then the 10 metric values in the same order as
# the header row above. "-" means a non-numeric placeholder value.
$rows = @(
  @{
    Period = "set_2024_a_ago_2025"
    Values = @(109318007711.96001, -844405.5, 12049575.5, $dash, $dash,
               109306802541.96001, 53034127686.489998, 53560333245.559998,
               50882316583.279999, 48204299921)
  },
  @{
    Period = "mai_2025_a_abr_2025"
    Values = @(106688331960.66, 19069999.73, 11634838.5, $dash, $dash,
               106657627122.42999, 52056327789.269997, 52262237289.989998,
               49649125425.489998, 47036013560.989998)
  },
  @{
    Period = "jan_2024_a_dez_2024"
    Values = @(103495630628.38, 18269291.93, 11634838.5, $dash, $dash,
               103465726497.95, 50492145171.779999, 50698205984,
               48163295684.800003, 45628385385.599998)
  },
  @{
    Period = "set_2023_a_ago_2024"
    Values = @(99192789740.149994, 56492032.229999997, 19022957, $dash, $dash,
               99117274750.919998, 49676872959.790001, 48567464627.949997,
               46139091396.550003, 43710718165.160004)
  },
  @{
    Period = "mai_2023_a_abr_2024"
    Values = @(95864466156.229996, 79853284.640000001, 18522957, $dash, $dash,
               95766089914.589996, 48234855928.120003, 46925384058.150002,
               44579114855.239998, 42232845652.330002)
  }
)

for ($r = 0; $r -lt $rows.Length; $r++) {
  $row = $r + 2   # row2 .. row6
  $ws.Cells.Item($row, 1).Value = $rows[$r].Period

  $values = $rows[$r].Values
  for ($c = 0; $c -lt $values.Length; $c++) {
    $col = $c + 2   # B=2 .. K=11
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $values[$c]
    if ($values[$c] -ne $dash) {
      $cell.NumberFormat = $numFmt
    }
  }
}

# ---- Cosmetic sheet-level tweaks to mirror the reshaped layout -----------
$ws.Columns.Item(1).ColumnWidth = 19.75
$ws.Range("A1:XFD1048576").Select() | Out-Null

